# Scheduled-runner refresh of Sheets/Cactuar_Profits.xlsx market data.
# Updates currentAveragePrice(NQ/HQ) and derived Leve-profit columns (H:N)
# per worksheet/row with freshly pulled marketboard values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 726.7059
$ws.Range("I80").Value = 505.58334
$ws.Range("K80").Value = 1516.75002
$ws.Range("M80").Value = -518.7500199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 726.7059
$ws.Range("I83").Value = 505.58334
$ws.Range("K83").Value = 4550.25006
$ws.Range("M83").Value = 441.7499399999997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1542019
$ws.Range("I86").Value = 2695939.8
$ws.Range("J86").Value = 3457.889
$ws.Range("K86").Value = 2695939.8
$ws.Range("L86").Value = 3457.889
$ws.Range("M86").Value = -2694816.8
$ws.Range("N86").Value = -5703.889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1542019
$ws.Range("I89").Value = 2695939.8
$ws.Range("J89").Value = 3457.889
$ws.Range("K89").Value = 13479699
$ws.Range("L89").Value = 17289.445
$ws.Range("M89").Value = -13474083
$ws.Range("N89").Value = -28521.445

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2531.0833
$ws.Range("I96").Value = 403
$ws.Range("K96").Value = 1209
$ws.Range("M96").Value = 164

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2597.2727
$ws.Range("I125").Value = 1068
$ws.Range("K125").Value = 9612
$ws.Range("M125").Value = -7152

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 116951.8
$ws.Range("J136").Value = 116951.8
$ws.Range("L136").Value = 116951.8
$ws.Range("N136").Value = -127151.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3143.1738
$ws.Range("I137").Value = 938.9259
$ws.Range("J137").Value = 6275.5264
$ws.Range("K137").Value = 2816.7777
$ws.Range("L137").Value = 18826.5792
$ws.Range("M137").Value = -266.7776999999996
$ws.Range("N137").Value = -23926.5792

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3533
$ws.Range("I138").Value = 1589.1428
$ws.Range("J138").Value = 4213.35
$ws.Range("K138").Value = 4767.428400000001
$ws.Range("L138").Value = 12640.05
$ws.Range("M138").Value = 372.5715999999993
$ws.Range("N138").Value = -22920.05

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 123999.75
$ws.Range("J139").Value = 123999.75
$ws.Range("L139").Value = 123999.75
$ws.Range("N139").Value = -134279.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 56401.637
$ws.Range("J140").Value = 54970.9
$ws.Range("L140").Value = 54970.9
$ws.Range("N140").Value = -65330.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1343603.1
$ws.Range("I2").Value = 2181931
$ws.Range("J2").Value = 2278.4
$ws.Range("K2").Value = 2181931
$ws.Range("L2").Value = 2278.4
$ws.Range("M2").Value = -2181818
$ws.Range("N2").Value = -2504.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3552.4473
$ws.Range("J32").Value = 12738
$ws.Range("L32").Value = 12738
$ws.Range("N32").Value = -13312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4019.2
$ws.Range("J63").Value = 3800.5
$ws.Range("L63").Value = 3800.5
$ws.Range("N63").Value = -5172.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4019.2
$ws.Range("J66").Value = 3800.5
$ws.Range("L66").Value = 19002.5
$ws.Range("N66").Value = -25866.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7814040.5
$ws.Range("I74").Value = 13158948
$ws.Range("K74").Value = 13158948
$ws.Range("M74").Value = -13158074

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7814040.5
$ws.Range("I77").Value = 13158948
$ws.Range("K77").Value = 65794740
$ws.Range("M77").Value = -65790372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 298899.94
$ws.Range("I102").Value = 508183.06
$ws.Range("K102").Value = 508183.06
$ws.Range("M102").Value = -506561.06

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1343603.1
$ws.Range("I116").Value = 2181931
$ws.Range("J116").Value = 2278.4
$ws.Range("K116").Value = 2181931
$ws.Range("L116").Value = 2278.4
$ws.Range("M116").Value = -2179637
$ws.Range("N116").Value = -6866.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2182.5789
$ws.Range("J122").Value = 1760
$ws.Range("L122").Value = 5280
$ws.Range("N122").Value = -10180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 55001
$ws.Range("I134").Value = 40001
$ws.Range("J134").Value = 70001
$ws.Range("K134").Value = 40001
$ws.Range("L134").Value = 70001
$ws.Range("M134").Value = -34931
$ws.Range("N134").Value = -80141

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1343603.1
$ws.Range("I3").Value = 2181931
$ws.Range("J3").Value = 2278.4
$ws.Range("K3").Value = 2181931
$ws.Range("L3").Value = 2278.4
$ws.Range("M3").Value = -2181817
$ws.Range("N3").Value = -2506.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2396.0588
$ws.Range("I20").Value = 2904.8333
$ws.Range("J20").Value = 1175
$ws.Range("K20").Value = 2904.8333
$ws.Range("L20").Value = 1175
$ws.Range("M20").Value = -2657.8333
$ws.Range("N20").Value = -1669

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1808.5238
$ws.Range("I16").Value = 1717.9231
$ws.Range("J16").Value = 1955.75
$ws.Range("K16").Value = 1717.9231
$ws.Range("L16").Value = 1955.75
$ws.Range("M16").Value = -1430.9231
$ws.Range("N16").Value = -2529.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3698.9583
$ws.Range("I31").Value = 1023.75
$ws.Range("J31").Value = 6374.1665
$ws.Range("K31").Value = 1023.75
$ws.Range("L31").Value = 6374.1665
$ws.Range("M31").Value = -728.75
$ws.Range("N31").Value = -6964.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3698.9583
$ws.Range("I34").Value = 1023.75
$ws.Range("J34").Value = 6374.1665
$ws.Range("K34").Value = 1023.75
$ws.Range("L34").Value = 6374.1665
$ws.Range("M34").Value = -821.75
$ws.Range("N34").Value = -6778.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 419305.12
$ws.Range("I58").Value = 668929.2
$ws.Range("K58").Value = 668929.2
$ws.Range("M58").Value = -668726.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1808.5238
$ws.Range("I113").Value = 1717.9231
$ws.Range("J113").Value = 1955.75
$ws.Range("K113").Value = 1717.9231
$ws.Range("L113").Value = 1955.75
$ws.Range("M113").Value = 452.0769
$ws.Range("N113").Value = -6295.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1917.52
$ws.Range("I122").Value = 1818.25
$ws.Range("J122").Value = 2094
$ws.Range("K122").Value = 5454.75
$ws.Range("L122").Value = 6282
$ws.Range("M122").Value = -3004.75
$ws.Range("N122").Value = -11182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10755008
$ws.Range("I132").Value = 13890971
$ws.Range("K132").Value = 41672913
$ws.Range("M132").Value = -41670383

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2579.7886
$ws.Range("I134").Value = 2697.2444
$ws.Range("K134").Value = 8091.733200000001
$ws.Range("M134").Value = -5556.733200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 419305.12
$ws.Range("I136").Value = 668929.2
$ws.Range("K136").Value = 2006787.6
$ws.Range("M136").Value = -2004237.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 92134.71000000001
$ws.Range("J141").Value = 100908.164
$ws.Range("L141").Value = 100908.164
$ws.Range("N141").Value = -111268.164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 4771.3335
$ws.Range("J52").Value = 4771.3335
$ws.Range("L52").Value = 14314.0005
$ws.Range("N52").Value = -14846.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 4704.6665
$ws.Range("J117").Value = 7000
$ws.Range("L117").Value = 21000
$ws.Range("N117").Value = -27884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1676.7778
$ws.Range("J129").Value = 1536.5
$ws.Range("L129").Value = 4609.5
$ws.Range("N129").Value = -14609.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 8477449
$ws.Range("I137").Value = 5000
$ws.Range("J137").Value = 9183486
$ws.Range("K137").Value = 15000
$ws.Range("L137").Value = 27550458
$ws.Range("M137").Value = -9900
$ws.Range("N137").Value = -27560658

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 75924.75
$ws.Range("J140").Value = 75924.75
$ws.Range("L140").Value = 75924.75
$ws.Range("N140").Value = -86284.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 278
$ws.Range("J55").Value = 263.66666
$ws.Range("L55").Value = 263.66666
$ws.Range("N55").Value = -609.66666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 627.6
$ws.Range("I113").Value = 673.4545000000001
$ws.Range("J113").Value = 501.5
$ws.Range("K113").Value = 2020.3635
$ws.Range("L113").Value = 1504.5
$ws.Range("M113").Value = 149.6364999999998
$ws.Range("N113").Value = -5844.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1048.9131
$ws.Range("I126").Value = 1040.5
$ws.Range("J126").Value = 1079.2
$ws.Range("K126").Value = 3121.5
$ws.Range("L126").Value = 3237.6
$ws.Range("M126").Value = -651.5
$ws.Range("N126").Value = -8177.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10486799
$ws.Range("I132").Value = 2417442.5
$ws.Range("K132").Value = 7252327.5
$ws.Range("M132").Value = -7249797.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 77000
$ws.Range("J138").Value = 77000
$ws.Range("L138").Value = 77000
$ws.Range("N138").Value = -87280
